{"js": "// The document defines three \"conditionally hide\" placeholders that used the\n// old `:collapse:hide` modifier syntax. The edit renames that modifier to the\n// newer `:hide-block-if-empty` syntax everywhere it appears:\n//   {{contacts}:collapse:hide}          -> {{contacts}:hide-block-if-empty}\n//   {{employees.hireDate}:collapse:hide} -> {{employees.hireDate}:hide-block-if-empty} (x2)\nconst body = context.document.body;\n\nconst hits = body.search(\"collapse:hide\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < hits.items.length; i++) {\n  hits.items[i].insertText(\"hide-block-if-empty\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The document defines three \"conditionally hide\" placeholders that used the\n# old `:collapse:hide` modifier syntax. The edit renames that modifier to the\n# newer `:hide-block-if-empty` syntax everywhere it appears:\n#   {{contacts}:collapse:hide}           -> {{contacts}:hide-block-if-empty}\n#   {{employees.hireDate}:collapse:hide} -> {{employees.hireDate}:hide-block-if-empty} (x2)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute(\"collapse:hide\", $false, $false, $false, $false, $false, $true, 1, $true, \"hide-block-if-empty\", 2)\n"}
